# Updated cryptos list with GitHub Actions
# Writes refreshed price/volume figures (and two name swaps caused by
# a ranking reshuffle: Chainlink<->Polygon at rows 17-18, FLOKI<->Stellar
# at rows 47-48) into the existing "Coin / Link / Price / Volume(1h)" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force text storage so numeric-looking strings (prices like
    # "20.86" or "0.999") are not reinterpreted as numbers, and restore
    # the default cell style afterwards so no stray formatting is left behind.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "73.306.11"
Set-TextCell 2 5 "  -0.24%  "

# Row 3
Set-TextCell 3 4 "3.973.61"
Set-TextCell 3 5 "  -2.21%  "

# Row 4
Set-TextCell 4 5 "  +0.04%  "

# Row 5
Set-TextCell 5 4 "608.42"
Set-TextCell 5 5 "  +5.65%  "

# Row 6
Set-TextCell 6 4 "169.37"
Set-TextCell 6 5 "  +11.05%  "

# Row 8
Set-TextCell 8 4 "0.999"
Set-TextCell 8 5 "  +0.01%  "

# Row 9
Set-TextCell 9 4 "0.788"
Set-TextCell 9 5 "  +2.46%  "

# Row 10
Set-TextCell 10 4 "0.185"
Set-TextCell 10 5 "  +7.50%  "

# Row 11
Set-TextCell 11 4 "56.57"
Set-TextCell 11 5 "  +4.02%  "

# Row 12
Set-TextCell 12 4 "0.0000335"
Set-TextCell 12 5 "  +1.46%  "

# Row 13
Set-TextCell 13 4 "11.32"
Set-TextCell 13 5 "  +1.01%  "

# Row 14
Set-TextCell 14 4 "4.609.79"
Set-TextCell 14 5 "  -2.22%  "

# Row 15
Set-TextCell 15 4 "3.984.40"
Set-TextCell 15 5 "  -2.06%  "

# Row 16
Set-TextCell 16 4 "14.25"
Set-TextCell 16 5 "  -1.57%  "

# Row 17
Set-TextCell 17 2 "Chainlink"
Set-TextCell 17 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 17 4 "20.86"
Set-TextCell 17 5 "  +0.02%  "

# Row 18
Set-TextCell 18 2 "Polygon"
Set-TextCell 18 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 18 4 "1.24"
Set-TextCell 18 5 "  +1.09%  "

# Row 19
Set-TextCell 19 4 "73.233.48"
Set-TextCell 19 5 "  -0.13%  "

# Row 20
Set-TextCell 20 5 "  -1.31%  "

# Row 21
Set-TextCell 21 4 "459.52"
Set-TextCell 21 5 "  +2.95%  "

# Row 22
Set-TextCell 22 4 "4.84"
Set-TextCell 22 5 "  +4.61%  "

# Row 23
Set-TextCell 23 4 "96.07"
Set-TextCell 23 5 "  -2.51%  "

# Row 24
Set-TextCell 24 4 "3.38"
Set-TextCell 24 5 "  -5.97%  "

# Row 25
Set-TextCell 25 4 "14.22"
Set-TextCell 25 5 "  -4.07%  "

# Row 26
Set-TextCell 26 4 "4.19"
Set-TextCell 26 5 "  -2.12%  "

# Row 27
Set-TextCell 27 5 "  -3.18%  "

# Row 28
Set-TextCell 28 4 "5.96"
Set-TextCell 28 5 "  +0.00%  "

# Row 29
Set-TextCell 29 4 "10.52"
Set-TextCell 29 5 "  -5.26%  "

# Row 30
Set-TextCell 30 4 "36.28"
Set-TextCell 30 5 "  -2.70%  "

# Row 31
Set-TextCell 31 4 "7.98"
Set-TextCell 31 5 "  +0.93%  "

# Row 32
Set-TextCell 32 4 "13.90"
Set-TextCell 32 5 "  +1.83%  "

# Row 33
Set-TextCell 33 4 "0.0000105"
Set-TextCell 33 5 "  +15.64%  "

# Row 34
Set-TextCell 34 4 "0.129"
Set-TextCell 34 5 "  -3.82%  "

# Row 35
Set-TextCell 35 4 "48.06"
Set-TextCell 35 5 "  -1.43%  "

# Row 36
Set-TextCell 36 4 "70.30"
Set-TextCell 36 5 "  +2.99%  "

# Row 37
Set-TextCell 37 4 "639.07"
Set-TextCell 37 5 "  -7.13%  "

# Row 38
Set-TextCell 38 4 "0.430"
Set-TextCell 38 5 "  -3.87%  "

# Row 39
Set-TextCell 39 5 "  -1.17%  "

# Row 40
Set-TextCell 40 4 "3.39"
Set-TextCell 40 5 "  -0.14%  "

# Row 41
Set-TextCell 41 4 "0.999"
Set-TextCell 41 5 "  -0.03%  "

# Row 42
Set-TextCell 42 5 "  +0.16%  "

# Row 43
Set-TextCell 43 4 "3.25"
Set-TextCell 43 5 "  +40.33%  "

# Row 44
Set-TextCell 44 4 "0.0483"
Set-TextCell 44 5 "  -3.50%  "

# Row 45
Set-TextCell 45 4 "10.62"
Set-TextCell 45 5 "  -5.91%  "

# Row 46
Set-TextCell 46 4 "3.15"
Set-TextCell 46 5 "  -5.76%  "

# Row 47
Set-TextCell 47 2 "FLOKI"
Set-TextCell 47 3 "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextCell 47 4 "0.000304"
Set-TextCell 47 5 "  +9.29%  "

# Row 48
Set-TextCell 48 2 "Stellar"
Set-TextCell 48 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 48 4 "0.149"
Set-TextCell 48 5 "  -3.01%  "

# Row 49
Set-TextCell 49 4 "3.45"
Set-TextCell 49 5 "  +3.25%  "

# Row 50
Set-TextCell 50 4 "2.58"
Set-TextCell 50 5 "  -4.81%  "

# Row 51
Set-TextCell 51 5 "  -3.21%  "
